$d = $word.ActiveDocument

# Update the date heading
$d.Content.Find.Execute("2026-02-17 Tuesday", $true, $false, $false, $false, $false, $true, 1, $false, "2026-02-18 Wednesday", 2)

# Update the division-problem table, addressing cells by (row, column)
# so that values which coincide between an "old" cell and a "new" cell
# elsewhere in the table never collide during replacement.
$tbl = $d.Tables.Item(1)

$tbl.Cell(1,1).Range.Text  = "84÷4=21, 0"
$tbl.Cell(1,2).Range.Text  = "40÷8=5, 0"
$tbl.Cell(1,3).Range.Text  = "81÷8=10, 1"
$tbl.Cell(1,4).Range.Text  = "52÷2=26, 0"
$tbl.Cell(1,5).Range.Text  = "44÷6=7, 2"

$tbl.Cell(5,1).Range.Text  = "43÷7=6, 1"
$tbl.Cell(5,2).Range.Text  = "70÷5=14, 0"
$tbl.Cell(5,3).Range.Text  = "38÷2=19, 0"
$tbl.Cell(5,4).Range.Text  = "55÷4=13, 3"
$tbl.Cell(5,5).Range.Text  = "48÷5=9, 3"

$tbl.Cell(9,1).Range.Text  = "39÷9=4, 3"
$tbl.Cell(9,2).Range.Text  = "67÷6=11, 1"
$tbl.Cell(9,3).Range.Text  = "93÷6=15, 3"
$tbl.Cell(9,4).Range.Text  = "78÷4=19, 2"
$tbl.Cell(9,5).Range.Text  = "28÷3=9, 1"

$tbl.Cell(13,1).Range.Text = "78÷3=26, 0"
$tbl.Cell(13,2).Range.Text = "22÷5=4, 2"
$tbl.Cell(13,3).Range.Text = "12÷8=1, 4"
$tbl.Cell(13,4).Range.Text = "31÷5=6, 1"
$tbl.Cell(13,5).Range.Text = "54÷2=27, 0"

$tbl.Cell(17,1).Range.Text = "15÷3=5, 0"
$tbl.Cell(17,2).Range.Text = "97÷9=10, 7"
$tbl.Cell(17,3).Range.Text = "24÷5=4, 4"
$tbl.Cell(17,4).Range.Text = "20÷3=6, 2"
$tbl.Cell(17,5).Range.Text = "93÷2=46, 1"
